# Finish TestDeleteClue and TestCRM suites
#
# 1. PageElement sheet: add four new Name/Value rows (28-31) describing the
#    new "delete clue" page elements.
# 2. TestData sheet: point CompanyName/CustomerName at the new test fixture
#    ("贵阳回煞酒厂" / "HuiShaaCEO") and add the VerifyDeleteClueResult
#    verification row (11).

$wb = $excel.ActiveWorkbook

$wsPage = $wb.Worksheets.Item("PageElement")
$wsData = $wb.Worksheets.Item("TestData")

# --- PageElement: new rows for the clue-selection / batch-delete flow ---
$wsPage.Range("A28").Value = "ViewCluePage_ClueCheckBox1"
$wsPage.Range("B28").Value = "//span[text()='"

# Leading single-quote is Excel's "quote prefix" marker, so an extra
# leading quote is needed to have the literal quote preserved in the value.
$wsPage.Range("A29").Value = "ViewCluePage_ClueCheckBox2"
$wsPage.Range("B29").Value = "'']/../../..//input[@class='check_list']"

$wsPage.Range("A30").Value = "ViewCluePage_ClueSelection"
$wsPage.Range("B30").Value = "//a[contains(text(),'批量操作')]"

$wsPage.Range("A31").Value = "ViewCluePage_ClueDelete"
$wsPage.Range("B31").Value = "//a[@id='delete']"

# --- TestData: swap in the new company/customer fixture values ---
$wsData.Range("B5").Value = "贵阳回煞酒厂"
$wsData.Range("B10").Value = "HuiShaaCEO"

# New verification row for the delete-clue result
$wsData.Range("D11").Value = "VerifyDeleteClueResult"
$wsData.Range("E11").Value = "删除成功!"

# --- Restore view/selection state ---
# TestData's selection moved (to a cell past the used range); touch it
# without leaving TestData as the active sheet.
$wsData.Range("G15").Select()

# PageElement stays the active tab, selection on the last new row.
$wsPage.Activate()
$wsPage.Range("B31").Select()
